$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting rows 30:55 down to 31:56
$ws.Rows("30").Insert()

# Populate new row 30 with the latest week's data
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44484
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 100112024
$ws.Range("G30").Value = "Choclo"
$ws.Range("H30").Value = "Dulce o Americano"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 150
$ws.Range("K30").Value = 43000
$ws.Range("L30").Value = 45000
$ws.Range("M30").Value = 44067
$ws.Range("N30").Value = "$/malla 70 unidades"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 630
$ws.Range("Q30").Value = 70
$ws.Range("R30").Value = "Hortaliza"
